$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 98; existing rows 98-116 shift down to 100-118.
$ws.Rows("98:99").Insert()

# --- New row 98 ---
$ws.Range("A98").Value = 11
$ws.Range("B98").Value = "Vega Monumental Concepción"
$ws.Range("C98").Value = "Bíobío"
$ws.Range("D98").Value = 44754
$ws.Range("E98").Value = 8
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100102
$ws.Range("H98").Value = "Cítricos"
$ws.Range("I98").Value = 100102004
$ws.Range("J98").Value = "Mandarina"
$ws.Range("K98").Value = "Clementina"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 250
$ws.Range("N98").Value = 8000
$ws.Range("O98").Value = 8500
$ws.Range("P98").Value = 8300
$ws.Range("Q98").Value = "`$/caja 18 kilos"
$ws.Range("R98").Value = "Región de O'Higgins"
$ws.Range("S98").Value = 461
$ws.Range("T98").Value = 18

# --- New row 99 ---
$ws.Range("A99").Value = 11
$ws.Range("B99").Value = "Vega Monumental Concepción"
$ws.Range("C99").Value = "Bíobío"
$ws.Range("D99").Value = 44754
$ws.Range("E99").Value = 8
$ws.Range("F99").Value = "Fruta"
$ws.Range("G99").Value = 100102
$ws.Range("H99").Value = "Cítricos"
$ws.Range("I99").Value = 100102004
$ws.Range("J99").Value = "Mandarina"
$ws.Range("K99").Value = "Clementina"
$ws.Range("L99").Value = "Segunda"
$ws.Range("M99").Value = 220
$ws.Range("N99").Value = 6500
$ws.Range("O99").Value = 7000
$ws.Range("P99").Value = 6727
$ws.Range("Q99").Value = "`$/caja 18 kilos"
$ws.Range("R99").Value = "Región de O'Higgins"
$ws.Range("S99").Value = 374
$ws.Range("T99").Value = 18
